# Insert a new weekly price record as row 170 on the single data sheet.
# Every existing record at/after row 170 (original rows 170-255) shifts
# down by one (to 171-256), and the newly inserted row 170 is populated
# with a new "Ciboulette" price observation for Región Metropolitana.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 170..255 down to 171..256, leaving a blank row 170 behind.
$ws.Rows.Item(170).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(170, 1).Value  = 4
$ws.Cells.Item(170, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(170, 3).Value  = "Los Lagos"
$ws.Cells.Item(170, 4).Value  = 44806
$ws.Cells.Item(170, 5).Value  = 10
$ws.Cells.Item(170, 6).Value  = 100112039
$ws.Cells.Item(170, 7).Value  = "Ciboulette"
$ws.Cells.Item(170, 8).Value  = "Sin especificar"
$ws.Cells.Item(170, 9).Value  = "Primera"
$ws.Cells.Item(170, 10).Value = 240
$ws.Cells.Item(170, 11).Value = 3000
$ws.Cells.Item(170, 12).Value = 3500
$ws.Cells.Item(170, 13).Value = 3250
$ws.Cells.Item(170, 14).Value = "$/docena de atados"
$ws.Cells.Item(170, 15).Value = "Región Metropolitana"
$ws.Cells.Item(170, 16).Value = 1083
$ws.Cells.Item(170, 17).Value = 3
$ws.Cells.Item(170, 18).Value = "Hortaliza"
